$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single-line mailing address into two lines (street / city-state-zip)
#    and add a blank paragraph after it. Using the "^p" paragraph-mark code in the
#    replacement text (rather than assigning Range.Text with embedded "`r") lets the
#    new paragraphs inherit formatting cleanly without leaving a stray empty run
#    behind in the new blank paragraph.
$d.Content.Find.Execute("15090 PARK DRIVE, SARATOGA CA 95070", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "15090 PARK DRIVE^pSARATOGA, CA 95070^p", 2) | Out-Null

# 3. Remove the two blank paragraphs (No Spacing, then Title) that used to sit right
#    after the "Board of Directors" line, leaving just the single Title-styled blank
#    paragraph that follows them.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13) -eq "Townhomes at Nuevo Homeowners Association Board of Directors") {
        $p1 = $d.Paragraphs.Item($i + 1)
        $p1.Range.Delete()
        $p2 = $d.Paragraphs.Item($i + 1)
        $p2.Range.Delete()
        break
    }
}
